$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameter")

# 1. Update the item placeholder text in row 23 (bug fix: proper casing / renamed field
#    so the edit-value-type message maps correctly).
$ws.Range("Q23").Value = "{{item.DescripcionTexto}}"
$ws.Range("R23").Value = "{{item.DescripcionParrafo}}"
$ws.Range("P23").Value = "{{item.Opcion}}"
$ws.Range("O23").Value = "{{item.MedidaTiempoId}}"

# 2. Un-merge the paired header cells on row 21 so each column header stands alone.
$ws.Range("D21:E21").UnMerge()
$ws.Range("F21:G21").UnMerge()
$ws.Range("H21:I21").UnMerge()
$ws.Range("J21:K21").UnMerge()
$ws.Range("L21:M21").UnMerge()
$ws.Range("N21:O21").UnMerge()
$ws.Range("P21:Q21").UnMerge()

# 3. Re-align the now-unmerged header cells: keep them vertically centered but drop the
#    forced horizontal centering (back to general alignment).
$ws.Range("D21:Q21").HorizontalAlignment = 1
$ws.Range("D21:Q21").VerticalAlignment = -4108

# 4. Leave the cursor where the author left it when they saved.
$ws.Range("F25").Select()
